# Exchange Signoff and LimeLite Issue fixed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSheet")

# Add the two new data rows at the bottom of the sheet
$ws.Range("A28").Value = "exchangesignoffserviceordernumber"
$ws.Range("B28").Value = "ZH00756"
$ws.Range("A29").Value = "PEL_PDL_LimeLitePackage"
$ws.Range("B29").Value = "Auto_PEL-PEL_WithoutRoutes (PEL)"

# Widen column A to fit the new, longer labels (no longer relying on Excel's
# "best fit" flag - a fixed custom width is used instead)
$ws.Columns.Item(1).ColumnWidth = 35.83

# Update the active selection to reflect where the user ended up after typing
$ws.Range("B30").Select() | Out-Null
